$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 428.41666
$ws.Range("I5").Value = 405.55554
$ws.Range("K5").Value = 405.55554
$ws.Range("M5").Value = -290.55554
$ws.Range("H18").Value = 1033.3334
$ws.Range("I18").Value = 1033.3334
$ws.Range("K18").Value = 1033.3334
$ws.Range("M18").Value = -749.3334
$ws.Range("H111").Value = 2870.8333
$ws.Range("I111").Value = 1298.25
$ws.Range("K111").Value = 3894.75
$ws.Range("M111").Value = -827.75
$ws.Range("H132").Value = 1110.2858
$ws.Range("I132").Value = 1146.6154
$ws.Range("J132").Value = 638
$ws.Range("K132").Value = 3439.8462
$ws.Range("L132").Value = 1914
$ws.Range("M132").Value = -909.8462
$ws.Range("N132").Value = -6974
$ws.Range("H141").Value = 2547.5
$ws.Range("I141").Value = 2547.5
$ws.Range("K141").Value = 7642.5
$ws.Range("M141").Value = -2462.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2505
$ws.Range("I2").Value = 3063.2
$ws.Range("J2").Value = 1109.5
$ws.Range("K2").Value = 3063.2
$ws.Range("L2").Value = 1109.5
$ws.Range("M2").Value = -2950.2
$ws.Range("N2").Value = -1335.5
$ws.Range("H32").Value = 16050.167
$ws.Range("I32").Value = 16050.167
$ws.Range("K32").Value = 16050.167
$ws.Range("M32").Value = -15763.167
$ws.Range("H45").Value = 2492.5833
$ws.Range("I45").Value = 2623.818
$ws.Range("K45").Value = 2623.818
$ws.Range("M45").Value = -2246.818
$ws.Range("H97").Value = 261.25
$ws.Range("I97").Value = 265.33334
$ws.Range("K97").Value = 265.33334
$ws.Range("M97").Value = 230.66666
$ws.Range("H102").Value = 8554.571
$ws.Range("I102").Value = 9147
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 9147
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -7525
$ws.Range("N102").Value = -8244
$ws.Range("H116").Value = 2505
$ws.Range("I116").Value = 3063.2
$ws.Range("J116").Value = 1109.5
$ws.Range("K116").Value = 3063.2
$ws.Range("L116").Value = 1109.5
$ws.Range("M116").Value = -769.1999999999998
$ws.Range("N116").Value = -5697.5
$ws.Range("H132").Value = 2043.3158
$ws.Range("I132").Value = 1296.0769
$ws.Range("J132").Value = 3662.3333
$ws.Range("K132").Value = 3888.2307
$ws.Range("L132").Value = 10986.9999
$ws.Range("M132").Value = -1358.2307
$ws.Range("N132").Value = -16046.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2505
$ws.Range("I3").Value = 3063.2
$ws.Range("J3").Value = 1109.5
$ws.Range("K3").Value = 3063.2
$ws.Range("L3").Value = 1109.5
$ws.Range("M3").Value = -2949.2
$ws.Range("N3").Value = -1337.5
$ws.Range("H105").Value = 3664.3333
$ws.Range("I105").Value = 3597.2
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3597.2
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1850.2
$ws.Range("N105").Value = -7494
$ws.Range("H107").Value = 1191.4546
$ws.Range("I107").Value = 1300.7778
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 1300.7778
$ws.Range("L107").Value = 699.5
$ws.Range("M107").Value = 619.2221999999999
$ws.Range("N107").Value = -4539.5
$ws.Range("H134").Value = 2291.8572
$ws.Range("I134").Value = 2207.1667
$ws.Range("K134").Value = 6621.500100000001
$ws.Range("M134").Value = -4086.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 959.4
$ws.Range("I16").Value = 1049.75
$ws.Range("J16").Value = 598
$ws.Range("K16").Value = 1049.75
$ws.Range("L16").Value = 598
$ws.Range("M16").Value = -762.75
$ws.Range("N16").Value = -1172
$ws.Range("H31").Value = 1895.8667
$ws.Range("I31").Value = 1230.25
$ws.Range("J31").Value = 2656.5715
$ws.Range("K31").Value = 1230.25
$ws.Range("L31").Value = 2656.5715
$ws.Range("M31").Value = -935.25
$ws.Range("N31").Value = -3246.5715
$ws.Range("H34").Value = 1895.8667
$ws.Range("I34").Value = 1230.25
$ws.Range("J34").Value = 2656.5715
$ws.Range("K34").Value = 1230.25
$ws.Range("L34").Value = 2656.5715
$ws.Range("M34").Value = -1028.25
$ws.Range("N34").Value = -3060.5715
$ws.Range("H99").Value = 4506
$ws.Range("I99").Value = 4506
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4506
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3008
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1143.95
$ws.Range("I107").Value = 1464.25
$ws.Range("J107").Value = 663.5
$ws.Range("K107").Value = 1464.25
$ws.Range("L107").Value = 663.5
$ws.Range("M107").Value = 455.75
$ws.Range("N107").Value = -4503.5
$ws.Range("H113").Value = 959.4
$ws.Range("I113").Value = 1049.75
$ws.Range("J113").Value = 598
$ws.Range("K113").Value = 1049.75
$ws.Range("L113").Value = 598
$ws.Range("M113").Value = 1120.25
$ws.Range("N113").Value = -4938
$ws.Range("H126").Value = 4506
$ws.Range("I126").Value = 4506
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13518
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11048
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3215.913
$ws.Range("J122").Value = 3379.8572
$ws.Range("L122").Value = 30418.7148
$ws.Range("N122").Value = -35318.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 9900
$ws.Range("I97").Value = 9900
$ws.Range("K97").Value = 9900
$ws.Range("M97").Value = -9404
$ws.Range("H107").Value = 627.1429000000001
$ws.Range("I107").Value = 1037
$ws.Range("K107").Value = 1037
$ws.Range("M107").Value = 883

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3255.889
$ws.Range("I40").Value = 3324.75
$ws.Range("K40").Value = 3324.75
$ws.Range("M40").Value = -3188.75
$ws.Range("H122").Value = 7498.25
$ws.Range("I122").Value = 7498.25
$ws.Range("K122").Value = 22494.75
$ws.Range("M122").Value = -20044.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1827.625
$ws.Range("I81").Value = 1827.625
$ws.Range("K81").Value = 3655.25
$ws.Range("M81").Value = -2594.25
$ws.Range("H84").Value = 1827.625
$ws.Range("I84").Value = 1827.625
$ws.Range("K84").Value = 18276.25
$ws.Range("M84").Value = -12972.25
$ws.Range("H96").Value = 3099.0833
$ws.Range("I96").Value = 3773.5557
$ws.Range("J96").Value = 1075.6666
$ws.Range("K96").Value = 3773.5557
$ws.Range("L96").Value = 1075.6666
$ws.Range("M96").Value = -2400.5557
$ws.Range("N96").Value = -3821.6666
$ws.Range("H126").Value = 1629.7778
$ws.Range("I126").Value = 1617
$ws.Range("J126").Value = 1640
$ws.Range("K126").Value = 4851
$ws.Range("L126").Value = 4920
$ws.Range("M126").Value = -2381
$ws.Range("N126").Value = -9860
